# Weekly fruit/vegetable price update: insert a new weekly record as
# row 70 (pushing the existing rows 70-78 down to 71-79) and populate it
# with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 70; this shifts rows 70-78 down to
# 71-79 and naturally inherits row 71's (formerly row 70's) formatting,
# including the date-column number format style used in column D.
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new weekly record.
$ws.Range("A70").Value = 11
$ws.Range("B70").Value = "Vega Monumental Concepción"
$ws.Range("C70").Value = "Bíobío"
$ws.Range("D70").Value = 44491
$ws.Range("E70").Value = 8
$ws.Range("F70").Value = 100112043
$ws.Range("G70").Value = "Pepino ensalada"
$ws.Range("H70").Value = "Sin especificar"
$ws.Range("I70").Value = "Primera"
$ws.Range("J70").Value = 100
$ws.Range("K70").Value = 8500
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = 8750
$ws.Range("N70").Value = "$/caja 60 unidades"
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 146
$ws.Range("Q70").Value = 60
$ws.Range("R70").Value = "Hortaliza"
